$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 230, shifting the existing rows 230:340 down to 231:341.
$ws.Rows(230).Insert()

# Populate the newly inserted row 230 with the new record.
$ws.Cells.Item(230, 1).Value = 5
$ws.Cells.Item(230, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(230, 3).Value = "Maule"
$ws.Cells.Item(230, 4).Value = "2023-07-27"
$ws.Cells.Item(230, 5).Value = 7
$ws.Cells.Item(230, 6).Value = 100112017
$ws.Cells.Item(230, 7).Value = "Apio"
$ws.Cells.Item(230, 8).Value = "Americana (o)"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 700
$ws.Cells.Item(230, 11).Value = 5500
$ws.Cells.Item(230, 12).Value = 5500
$ws.Cells.Item(230, 13).Value = 5500
$ws.Cells.Item(230, 14).Value = "$/docena de matas"
$ws.Cells.Item(230, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(230, 16).Value = 917
$ws.Cells.Item(230, 17).Value = 6
$ws.Cells.Item(230, 18).Value = "Hortaliza"
